$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "86÷6=14, 2"
$t.Cell(1,2).Range.Text = "66÷6=11, 0"
$t.Cell(1,3).Range.Text = "89÷5=17, 4"
$t.Cell(1,4).Range.Text = "29÷3=9, 2"
$t.Cell(1,5).Range.Text = "55÷2=27, 1"
$t.Cell(5,1).Range.Text = "26÷3=8, 2"
$t.Cell(5,2).Range.Text = "85÷7=12, 1"
$t.Cell(5,3).Range.Text = "91÷2=45, 1"
$t.Cell(5,4).Range.Text = "84÷5=16, 4"
$t.Cell(5,5).Range.Text = "90÷5=18, 0"
$t.Cell(9,1).Range.Text = "33÷4=8, 1"
$t.Cell(9,2).Range.Text = "80÷4=20, 0"
$t.Cell(9,3).Range.Text = "90÷3=30, 0"
$t.Cell(9,4).Range.Text = "40÷2=20, 0"
$t.Cell(9,5).Range.Text = "52÷7=7, 3"
$t.Cell(13,1).Range.Text = "23÷4=5, 3"
$t.Cell(13,2).Range.Text = "10÷6=1, 4"
$t.Cell(13,3).Range.Text = "99÷7=14, 1"
$t.Cell(13,4).Range.Text = "24÷8=3, 0"
$t.Cell(13,5).Range.Text = "39÷5=7, 4"
$t.Cell(17,1).Range.Text = "97÷4=24, 1"
$t.Cell(17,2).Range.Text = "32÷5=6, 2"
$t.Cell(17,3).Range.Text = "95÷3=31, 2"
$t.Cell(17,4).Range.Text = "79÷4=19, 3"
$t.Cell(17,5).Range.Text = "66÷4=16, 2"
